$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-12-03 Tuesday", $true, $true, $false, $false, $false, $true, 1, $false, "2024-12-04 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("98-50=", $true, $true, $false, $false, $false, $true, 1, $false, "51-44=", 2) | Out-Null
$d.Content.Find.Execute("58+24=", $true, $true, $false, $false, $false, $true, 1, $false, "81-76=", 2) | Out-Null
$d.Content.Find.Execute("7+30=", $true, $true, $false, $false, $false, $true, 1, $false, "53-5=", 2) | Out-Null
$d.Content.Find.Execute("18+29=", $true, $true, $false, $false, $false, $true, 1, $false, "87-85=", 2) | Out-Null
$d.Content.Find.Execute("67-24=", $true, $true, $false, $false, $false, $true, 1, $false, "85-61=", 2) | Out-Null
$d.Content.Find.Execute("42+8=", $true, $true, $false, $false, $false, $true, 1, $false, "71+24=", 2) | Out-Null
$d.Content.Find.Execute("23+34=", $true, $true, $false, $false, $false, $true, 1, $false, "42-29=", 2) | Out-Null
$d.Content.Find.Execute("17+41=", $true, $true, $false, $false, $false, $true, 1, $false, "72+6=", 2) | Out-Null
$d.Content.Find.Execute("92-9=", $true, $true, $false, $false, $false, $true, 1, $false, "83+12=", 2) | Out-Null
$d.Content.Find.Execute("30+63=", $true, $true, $false, $false, $false, $true, 1, $false, "93-31=", 2) | Out-Null
$d.Content.Find.Execute("91-64=", $true, $true, $false, $false, $false, $true, 1, $false, "36+59=", 2) | Out-Null
$d.Content.Find.Execute("0+97=", $true, $true, $false, $false, $false, $true, 1, $false, "33-8=", 2) | Out-Null
$d.Content.Find.Execute("25-6=", $true, $true, $false, $false, $false, $true, 1, $false, "92-26=", 2) | Out-Null
$d.Content.Find.Execute("77-52=", $true, $true, $false, $false, $false, $true, 1, $false, "66+5=", 2) | Out-Null
$d.Content.Find.Execute("53-26=", $true, $true, $false, $false, $false, $true, 1, $false, "3+95=", 2) | Out-Null
$d.Content.Find.Execute("4+77=", $true, $true, $false, $false, $false, $true, 1, $false, "79-32=", 2) | Out-Null
$d.Content.Find.Execute("91-79=", $true, $true, $false, $false, $false, $true, 1, $false, "99-36=", 2) | Out-Null
$d.Content.Find.Execute("11+23=", $true, $true, $false, $false, $false, $true, 1, $false, "15-8=", 2) | Out-Null
$d.Content.Find.Execute("53-40=", $true, $true, $false, $false, $false, $true, 1, $false, "64+23=", 2) | Out-Null
$d.Content.Find.Execute("44+39=", $true, $true, $false, $false, $false, $true, 1, $false, "58-37=", 2) | Out-Null
$d.Content.Find.Execute("19-14=", $true, $true, $false, $false, $false, $true, 1, $false, "69+29=", 2) | Out-Null
$d.Content.Find.Execute("31+24=", $true, $true, $false, $false, $false, $true, 1, $false, "48+24=", 2) | Out-Null
$d.Content.Find.Execute("66-20=", $true, $true, $false, $false, $false, $true, 1, $false, "7+19=", 2) | Out-Null
$d.Content.Find.Execute("45+6=", $true, $true, $false, $false, $false, $true, 1, $false, "77-74=", 2) | Out-Null
$d.Content.Find.Execute("66+32=", $true, $true, $false, $false, $false, $true, 1, $false, "51+15=", 2) | Out-Null
$d.Content.Find.Execute("9+81=", $true, $true, $false, $false, $false, $true, 1, $false, "69-54=", 2) | Out-Null
$d.Content.Find.Execute("4+48=", $true, $true, $false, $false, $false, $true, 1, $false, "50+29=", 2) | Out-Null
$d.Content.Find.Execute("13+60=", $true, $true, $false, $false, $false, $true, 1, $false, "57+9=", 2) | Out-Null
$d.Content.Find.Execute("28+47=", $true, $true, $false, $false, $false, $true, 1, $false, "91-89=", 2) | Out-Null
$d.Content.Find.Execute("17+8=", $true, $true, $false, $false, $false, $true, 1, $false, "42+29=", 2) | Out-Null
$d.Content.Find.Execute("14-11=", $true, $true, $false, $false, $false, $true, 1, $false, "62-5=", 2) | Out-Null
$d.Content.Find.Execute("51-6=", $true, $true, $false, $false, $false, $true, 1, $false, "5+43=", 2) | Out-Null
$d.Content.Find.Execute("11+2=", $true, $true, $false, $false, $false, $true, 1, $false, "0+18=", 2) | Out-Null
$d.Content.Find.Execute("8+89=", $true, $true, $false, $false, $false, $true, 1, $false, "86-71=", 2) | Out-Null
$d.Content.Find.Execute("14+56=", $true, $true, $false, $false, $false, $true, 1, $false, "69-18=", 2) | Out-Null
$d.Content.Find.Execute("91+3=", $true, $true, $false, $false, $false, $true, 1, $false, "43-25=", 2) | Out-Null
$d.Content.Find.Execute("58-18=", $true, $true, $false, $false, $false, $true, 1, $false, "48+25=", 2) | Out-Null
$d.Content.Find.Execute("66-36=", $true, $true, $false, $false, $false, $true, 1, $false, "12-4=", 2) | Out-Null
$d.Content.Find.Execute("8+87=", $true, $true, $false, $false, $false, $true, 1, $false, "32-6=", 2) | Out-Null
$d.Content.Find.Execute("88-45=", $true, $true, $false, $false, $false, $true, 1, $false, "94-49=", 2) | Out-Null
$d.Content.Find.Execute("86-60=", $true, $true, $false, $false, $false, $true, 1, $false, "39+26=", 2) | Out-Null
$d.Content.Find.Execute("20+50=", $true, $true, $false, $false, $false, $true, 1, $false, "28+31=", 2) | Out-Null
$d.Content.Find.Execute("32+6=", $true, $true, $false, $false, $false, $true, 1, $false, "0+5=", 2) | Out-Null
$d.Content.Find.Execute("73+11=", $true, $true, $false, $false, $false, $true, 1, $false, "29+25=", 2) | Out-Null
$d.Content.Find.Execute("92-22=", $true, $true, $false, $false, $false, $true, 1, $false, "90-8=", 2) | Out-Null
$d.Content.Find.Execute("59-57=", $true, $true, $false, $false, $false, $true, 1, $false, "59+22=", 2) | Out-Null
$d.Content.Find.Execute("53-42=", $true, $true, $false, $false, $false, $true, 1, $false, "44+28=", 2) | Out-Null
$d.Content.Find.Execute("16+15=", $true, $true, $false, $false, $false, $true, 1, $false, "29+16=", 2) | Out-Null
$d.Content.Find.Execute("2+69=", $true, $true, $false, $false, $false, $true, 1, $false, "79-75=", 2) | Out-Null
$d.Content.Find.Execute("18+42=", $true, $true, $false, $false, $false, $true, 1, $false, "6+80=", 2) | Out-Null
$d.Content.Find.Execute("88-30=", $true, $true, $false, $false, $false, $true, 1, $false, "47+32=", 2) | Out-Null
$d.Content.Find.Execute("50-13=", $true, $true, $false, $false, $false, $true, 1, $false, "95-81=", 2) | Out-Null
$d.Content.Find.Execute("85-83=", $true, $true, $false, $false, $false, $true, 1, $false, "70-29=", 2) | Out-Null
$d.Content.Find.Execute("36+11=", $true, $true, $false, $false, $false, $true, 1, $false, "92-19=", 2) | Out-Null
$d.Content.Find.Execute("38+35=", $true, $true, $false, $false, $false, $true, 1, $false, "6+64=", 2) | Out-Null
$d.Content.Find.Execute("53+18=", $true, $true, $false, $false, $false, $true, 1, $false, "89+1=", 2) | Out-Null
$d.Content.Find.Execute("23+8=", $true, $true, $false, $false, $false, $true, 1, $false, "56-55=", 2) | Out-Null
$d.Content.Find.Execute("61+38=", $true, $true, $false, $false, $false, $true, 1, $false, "10+49=", 2) | Out-Null
$d.Content.Find.Execute("47+5=", $true, $true, $false, $false, $false, $true, 1, $false, "59-1=", 2) | Out-Null
$d.Content.Find.Execute("30-26=", $true, $true, $false, $false, $false, $true, 1, $false, "41+55=", 2) | Out-Null
$d.Content.Find.Execute("46-37=", $true, $true, $false, $false, $false, $true, 1, $false, "30-7=", 2) | Out-Null
$d.Content.Find.Execute("87-35=", $true, $true, $false, $false, $false, $true, 1, $false, "89-13=", 2) | Out-Null
$d.Content.Find.Execute("92-77=", $true, $true, $false, $false, $false, $true, 1, $false, "14+35=", 2) | Out-Null
$d.Content.Find.Execute("29+64=", $true, $true, $false, $false, $false, $true, 1, $false, "93-6=", 2) | Out-Null
$d.Content.Find.Execute("17-11=", $true, $true, $false, $false, $false, $true, 1, $false, "40+49=", 2) | Out-Null
$d.Content.Find.Execute("10+21=", $true, $true, $false, $false, $false, $true, 1, $false, "85-77=", 2) | Out-Null
$d.Content.Find.Execute("15+80=", $true, $true, $false, $false, $false, $true, 1, $false, "70-64=", 2) | Out-Null
$d.Content.Find.Execute("76+0=", $true, $true, $false, $false, $false, $true, 1, $false, "64+19=", 2) | Out-Null
$d.Content.Find.Execute("54+10=", $true, $true, $false, $false, $false, $true, 1, $false, "86-4=", 2) | Out-Null
$d.Content.Find.Execute("33+13=", $true, $true, $false, $false, $false, $true, 1, $false, "26+49=", 2) | Out-Null
$d.Content.Find.Execute("86-55=", $true, $true, $false, $false, $false, $true, 1, $false, "28+56=", 2) | Out-Null
$d.Content.Find.Execute("12+86=", $true, $true, $false, $false, $false, $true, 1, $false, "99-55=", 2) | Out-Null
$d.Content.Find.Execute("87-87=", $true, $true, $false, $false, $false, $true, 1, $false, "92-49=", 2) | Out-Null
$d.Content.Find.Execute("39+50=", $true, $true, $false, $false, $false, $true, 1, $false, "63-55=", 2) | Out-Null
$d.Content.Find.Execute("4-4=", $true, $true, $false, $false, $false, $true, 1, $false, "81-63=", 2) | Out-Null
$d.Content.Find.Execute("54-33=", $true, $true, $false, $false, $false, $true, 1, $false, "80-19=", 2) | Out-Null
$d.Content.Find.Execute("89-74=", $true, $true, $false, $false, $false, $true, 1, $false, "1+25=", 2) | Out-Null
$d.Content.Find.Execute("98-98=", $true, $true, $false, $false, $false, $true, 1, $false, "90-81=", 2) | Out-Null
$d.Content.Find.Execute("71+13=", $true, $true, $false, $false, $false, $true, 1, $false, "7+83=", 2) | Out-Null
$d.Content.Find.Execute("77+7=", $true, $true, $false, $false, $false, $true, 1, $false, "27+45=", 2) | Out-Null
$d.Content.Find.Execute("25-17=", $true, $true, $false, $false, $false, $true, 1, $false, "83-79=", 2) | Out-Null
$d.Content.Find.Execute("71-21=", $true, $true, $false, $false, $false, $true, 1, $false, "15+48=", 2) | Out-Null
$d.Content.Find.Execute("25+3=", $true, $true, $false, $false, $false, $true, 1, $false, "77-56=", 2) | Out-Null
$d.Content.Find.Execute("30-13=", $true, $true, $false, $false, $false, $true, 1, $false, "73-54=", 2) | Out-Null
$d.Content.Find.Execute("2+82=", $true, $true, $false, $false, $false, $true, 1, $false, "24-13=", 2) | Out-Null
$d.Content.Find.Execute("17+60=", $true, $true, $false, $false, $false, $true, 1, $false, "31+37=", 2) | Out-Null
$d.Content.Find.Execute("69-49=", $true, $true, $false, $false, $false, $true, 1, $false, "27-0=", 2) | Out-Null
$d.Content.Find.Execute("41-28=", $true, $true, $false, $false, $false, $true, 1, $false, "48-34=", 2) | Out-Null
$d.Content.Find.Execute("29-10=", $true, $true, $false, $false, $false, $true, 1, $false, "59-47=", 2) | Out-Null
$d.Content.Find.Execute("75-16=", $true, $true, $false, $false, $false, $true, 1, $false, "80-49=", 2) | Out-Null
$d.Content.Find.Execute("97-0=", $true, $true, $false, $false, $false, $true, 1, $false, "45-43=", 2) | Out-Null
$d.Content.Find.Execute("35+55=", $true, $true, $false, $false, $false, $true, 1, $false, "29-7=", 2) | Out-Null
$d.Content.Find.Execute("66+18=", $true, $true, $false, $false, $false, $true, 1, $false, "88+9=", 2) | Out-Null
$d.Content.Find.Execute("86-24=", $true, $true, $false, $false, $false, $true, 1, $false, "87-12=", 2) | Out-Null
$d.Content.Find.Execute("30+52=", $true, $true, $false, $false, $false, $true, 1, $false, "0+73=", 2) | Out-Null
$d.Content.Find.Execute("86-17=", $true, $true, $false, $false, $false, $true, 1, $false, "42+54=", 2) | Out-Null
$d.Content.Find.Execute("36-6=", $true, $true, $false, $false, $false, $true, 1, $false, "67+4=", 2) | Out-Null
$d.Content.Find.Execute("93-89=", $true, $true, $false, $false, $false, $true, 1, $false, "81-35=", 2) | Out-Null
$d.Content.Find.Execute("42+26=", $true, $true, $false, $false, $false, $true, 1, $false, "98-32=", 2) | Out-Null
$d.Content.Find.Execute("13+7=", $true, $true, $false, $false, $false, $true, 1, $false, "57+29=", 2) | Out-Null
